$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "60.973.96"
$ws.Cells.Item(2, 5).Value = "  -0.53%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.363.05"
$ws.Cells.Item(3, 5).Value = "  -3.78%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.10%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'541.48"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -1.17%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'137.93"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  -5.65%  "

# Row 7
$ws.Cells.Item(7, 5).Value = "  -0.04%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -11.28%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "2.359.66"
$ws.Cells.Item(9, 5).Value = "  -3.78%  "

# Row 10
$ws.Cells.Item(10, 5).Value = "  -1.50%  "

# Row 12
$ws.Cells.Item(12, 5).Value = "  -3.27%  "

# Row 13
$ws.Cells.Item(13, 4).Value = "'0.342"
$ws.Cells.Item(13, 4).Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -2.58%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "'24.91"
$ws.Cells.Item(14, 4).Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  -3.93%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "2.785.21"
$ws.Cells.Item(15, 5).Value = "  -3.96%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  -2.50%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "60.486.82"
$ws.Cells.Item(17, 5).Value = "  -1.21%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "2.359.85"
$ws.Cells.Item(18, 5).Value = "  -4.16%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "'10.63"
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  -4.02%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  -1.90%  "

# Row 21
$ws.Cells.Item(21, 4).Value = "'315.86"
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = "  -0.45%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'6.60"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -6.18%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'1.00"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  -0.04%  "

# Row 24
$ws.Cells.Item(24, 4).Value = "'1.91"
$ws.Cells.Item(24, 4).Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  +1.91%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'63.28"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -1.01%  "

# Row 26
$ws.Cells.Item(26, 4).Value = "'8.44"
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +10.69%  "

# Row 27
$ws.Cells.Item(27, 5).Value = "  +0.03%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "2.477.22"
$ws.Cells.Item(28, 5).Value = "  -4.07%  "

# Row 29
$ws.Cells.Item(29, 4).Value = "0.0₃0898"
$ws.Cells.Item(29, 5).Value = "  -7.11%  "

# Row 30
$ws.Cells.Item(30, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(30, 4).Value = "'7.96"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  -3.42%  "

# Row 31
$ws.Cells.Item(31, 2).Value = "Bittensor"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Cells.Item(31, 4).Value = "'508.41"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  -7.22%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  -4.93%  "

# Row 33
$ws.Cells.Item(33, 5).Value = "  -0.08%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -5.28%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'1.55"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -1.95%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  -0.12%  "

# Row 37
$ws.Cells.Item(37, 4).Value = "'4.61"
$ws.Cells.Item(37, 4).Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  -3.92%  "

# Row 38
$ws.Cells.Item(38, 4).Value = "'18.49"
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  +0.63%  "

# Row 39
$ws.Cells.Item(39, 5).Value = "  -1.49%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'5.26"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -10.13%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'1.79"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +2.15%  "

# Row 42
$ws.Cells.Item(42, 2).Value = "USDe"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(42, 4).Value = "'1.00"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  -0.01%  "

# Row 43
$ws.Cells.Item(43, 2).Value = "Monero"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(43, 4).Value = "'138.15"
$ws.Cells.Item(43, 4).Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  -2.76%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "'40.10"
$ws.Cells.Item(44, 4).Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  -0.91%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'2.14"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -8.73%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'138.75"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  -4.87%  "

# Row 47
$ws.Cells.Item(47, 5).Value = "  -1.49%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'0.0512"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -4.19%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'19.63"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -8.20%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "'0.572"
$ws.Cells.Item(50, 4).Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -2.63%  "

# Row 51
$ws.Cells.Item(51, 4).Value = "'0.0224"
$ws.Cells.Item(51, 4).Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -1.40%  "
